# Update column F (dSF) values for specific rows as part of a data repull /
# recalculation pass (see commit message: "repull data, push all data, mean
# calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -3
    7  = 5
    10 = 2
    13 = -2
    16 = -3
    18 = -7
    19 = -8
    23 = -3
    25 = -4
    26 = 1
    27 = -6
    30 = -4
    32 = -2
    34 = 5
    36 = -7
    37 = -5
    38 = 3
    42 = 0
    44 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
